$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MFG")

$ws.Range("D8").Value = 8339000
$ws.Range("E8").Value = 7588400
$ws.Range("F8").Value = 7114100
$ws.Range("G8").Value = 6527500
$ws.Range("H8").Value = 6929700
$ws.Range("I8").Value = 6631800
$ws.Range("J8").Value = 3556100

$ws.Range("D17").Value = 3997300
$ws.Range("E17").Value = 2632000
$ws.Range("F17").Value = 3592500
$ws.Range("G17").Value = 2363400
$ws.Range("H17").Value = 2638400
$ws.Range("I17").Value = 2003800
$ws.Range("J17").Value = 816000

$ws.Range("D18").Value = 4341800
$ws.Range("E18").Value = 4956400
$ws.Range("F18").Value = 3521600
$ws.Range("G18").Value = 4164100
$ws.Range("H18").Value = 4291300
$ws.Range("I18").Value = 4628000
$ws.Range("J18").Value = 2740100

$ws.Range("D20").Value = -1390800
$ws.Range("E20").Value = -320100
$ws.Range("F20").Value = -3312800
$ws.Range("G20").Value = -30300
$ws.Range("H20").Value = 1621500
$ws.Range("I20").Value = 276500
$ws.Range("J20").Value = -1057700

$ws.Range("D21").Value = 3756400
$ws.Range("E21").Value = 5450500
$ws.Range("F21").Value = 1018400
$ws.Range("G21").Value = 4871900
$ws.Range("H21").Value = 6721200

$ws.Range("D23").Value = 2951000
$ws.Range("E23").Value = 4636300
$ws.Range("F23").Value = 208800
$ws.Range("G23").Value = 4133700
$ws.Range("H23").Value = 5912800
$ws.Range("I23").Value = 4904500
$ws.Range("J23").Value = 1682500

$ws.Range("D24").Value = 1069200
$ws.Range("E24").Value = 1078800
$ws.Range("F24").Value = 151200
$ws.Range("G24").Value = 673600
$ws.Range("H24").Value = 1620700
$ws.Range("I24").Value = 1512000
$ws.Range("J24").Value = 697700

$ws.Range("D26").Value = 1881800
$ws.Range("E26").Value = 3557500
$ws.Range("F26").Value = 57600
$ws.Range("G26").Value = 3460100
$ws.Range("H26").Value = 4292100
$ws.Range("I26").Value = 3392400
$ws.Range("J26").Value = 984800

$ws.Range("D27").Value = 1853700
$ws.Range("E27").Value = 3367800
$ws.Range("F27").Value = -154700
$ws.Range("G27").Value = 3431200
$ws.Range("H27").Value = 4359000
$ws.Range("I27").Value = 3307500
$ws.Range("J27").Value = 782700

$ws.Range("D32").Value = 1390800
$ws.Range("E32").Value = 320100
$ws.Range("F32").Value = 3312800
$ws.Range("G32").Value = 30300
$ws.Range("H32").Value = -1621500
$ws.Range("I32").Value = -276500
$ws.Range("J32").Value = 1057700

$ws.Range("D33").Value = 1853700
$ws.Range("E33").Value = 3367800
$ws.Range("F33").Value = -154700
$ws.Range("G33").Value = 3431200
$ws.Range("H33").Value = 4359000
$ws.Range("I33").Value = 3307500
$ws.Range("J33").Value = 782700

$ws.Range("D35").Value = 1853700
$ws.Range("E35").Value = 3367800
$ws.Range("F35").Value = -154700
$ws.Range("G35").Value = 3431200
$ws.Range("H35").Value = 4359000
$ws.Range("I35").Value = 3307500
$ws.Range("J35").Value = 782700

$ws.Range("D41").Value = 435465100
$ws.Range("E41").Value = 464633800
$ws.Range("F41").Value = 430186200
$ws.Range("G41").Value = 391013700
$ws.Range("H41").Value = 331316000
$ws.Range("I41").Value = 320801800
$ws.Range("J41").Value = 263029300

$ws.Range("D42").Value = 359967800
$ws.Range("E42").Value = 400040300
$ws.Range("F42").Value = 352101800
$ws.Range("G42").Value = 406244600
$ws.Range("H42").Value = 386376500
$ws.Range("I42").Value = 382083700
$ws.Range("J42").Value = 291768200

$ws.Range("J47").Value = 2639400

$ws.Range("D48").Value = 19130300
$ws.Range("E48").Value = 18875000
$ws.Range("F48").Value = 18453100
$ws.Range("G48").Value = 17015300
$ws.Range("H48").Value = 16615400
$ws.Range("I48").Value = 15706700
$ws.Range("J48").Value = 9745600

$ws.Range("D49").Value = 1623900
$ws.Range("E49").Value = 1668300
$ws.Range("F49").Value = 1711500
$ws.Range("G49").Value = 588000
$ws.Range("H49").Value = 612400
$ws.Range("I49").Value = 558200
$ws.Range("J49").Value = 5944300

$ws.Range("D52").Value = 516100
$ws.Range("E52").Value = 562600
$ws.Range("F52").Value = 574200
$ws.Range("G52").Value = 842400
$ws.Range("H52").Value = 518400
$ws.Range("I52").Value = 513700
$ws.Range("J52").Value = 7054100

$ws.Range("D54").Value = 1846471000
$ws.Range("E54").Value = 1889152600
$ws.Range("F54").Value = 1812125000
$ws.Range("G54").Value = 1771973100
$ws.Range("H54").Value = 1752043800
$ws.Range("I54").Value = 1738731000
$ws.Range("J54").Value = 1714750100

$ws.Range("D57").Value = 16873800
$ws.Range("E57").Value = 28629300

$ws.Range("D59").Value = 2688400
$ws.Range("E59").Value = 2523700
$ws.Range("F59").Value = 2557500
$ws.Range("G59").Value = 2250000
$ws.Range("H59").Value = 2514500
$ws.Range("I59").Value = 2376700
$ws.Range("J59").Value = 541200

$ws.Range("D61").Value = 117115300
$ws.Range("E61").Value = 127758600
$ws.Range("F61").Value = 131345900
$ws.Range("G61").Value = 137775900
$ws.Range("H61").Value = 133480400
$ws.Range("I61").Value = 131837000
$ws.Range("J61").Value = 119414800

$ws.Range("D62").Value = 2768100
$ws.Range("E62").Value = 2131300
$ws.Range("F62").Value = 1269900
$ws.Range("G62").Value = 1368700
$ws.Range("H62").Value = 1824800
$ws.Range("I62").Value = 1810700
$ws.Range("J62").Value = 6508400

$ws.Range("D66").Value = 1766300500
$ws.Range("E66").Value = 1810490200
$ws.Range("F66").Value = 1737442300
$ws.Range("G66").Value = 1698720300
$ws.Range("H66").Value = 1679592200
$ws.Range("I66").Value = 1667324600
$ws.Range("J66").Value = 1640939000

$ws.Range("H70").Value = 894300
$ws.Range("I70").Value = 1302400
$ws.Range("J70").Value = 1926600

$ws.Range("D72").Value = 11807500
$ws.Range("E72").Value = 10814300
$ws.Range("F72").Value = 8306800
$ws.Range("G72").Value = 9321900
$ws.Range("H72").Value = 6750900
$ws.Range("I72").Value = 3206700
$ws.Range("J72").Value = 25069600

$ws.Range("D76").Value = 80170500
$ws.Range("E76").Value = 78662500
$ws.Range("F76").Value = 74682700
$ws.Range("G76").Value = 73252800
$ws.Range("H76").Value = 71557300
$ws.Range("I76").Value = 70104000
$ws.Range("J76").Value = 71884500

$ws.Range("D81").Value = 1853700
$ws.Range("E81").Value = 3367800
$ws.Range("F81").Value = -154700
$ws.Range("G81").Value = 3431200
$ws.Range("H81").Value = 4359000
$ws.Range("I81").Value = 3307500
$ws.Range("J81").Value = 782700

$ws.Range("D83").Value = 805400
$ws.Range("E83").Value = 814200
$ws.Range("F83").Value = 809600
$ws.Range("G83").Value = 738100
$ws.Range("H83").Value = 808400

$ws.Range("D89").Value = 21369300
$ws.Range("E89").Value = -22794500
$ws.Range("F89").Value = -5835000
$ws.Range("G89").Value = 16803300
$ws.Range("H89").Value = -36727700

$ws.Range("D91").Value = -1219300
$ws.Range("E91").Value = -1422200
$ws.Range("F91").Value = -1955400
$ws.Range("G91").Value = -1940000
$ws.Range("H91").Value = -3378600

$ws.Range("D94").Value = -560700
$ws.Range("E94").Value = -39102100
$ws.Range("F94").Value = -52396900
$ws.Range("G94").Value = -63329300
$ws.Range("H94").Value = -45693200

$ws.Range("D96").Value = -861000
$ws.Range("E96").Value = -860100
$ws.Range("F96").Value = -861100
$ws.Range("G96").Value = -856800
$ws.Range("H96").Value = -968500

$ws.Range("D100").Value = -18049600
$ws.Range("E100").Value = 60109300
$ws.Range("F100").Value = 59454100
$ws.Range("G100").Value = 47848500
$ws.Range("H100").Value = 56452900

$ws.Range("D101").Value = -132200
$ws.Range("E101").Value = 7100
$ws.Range("F101").Value = 438000
$ws.Range("G101").Value = -545600
$ws.Range("H101").Value = -98600

$ws.Range("D102").Value = 2626700
$ws.Range("E102").Value = -1780200
$ws.Range("F102").Value = 1660200
$ws.Range("G102").Value = 776900
$ws.Range("H102").Value = -26066600
